$d = $word.ActiveDocument

# 1) Merge the two adjacent runs " is clicked you should create" + " a "
#    (identical formatting) into a single run " is clicked you should create a ".
$rng = $d.Content
[void]$rng.Find.Execute(" is clicked you should create a ", $true, $false, $false, $false, $false, $true, 1, $false, " is clicked you should create a ", 2)

# 2) Move the "_GoBack" bookmark: remove it from its old location (end of the
#    paragraph containing the `"{author}: {message}"` code sample) ...
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}

# ... and re-create it right after the paragraph that ends with
# "to which your program should make requests is:" (immediately before the
# closing </w:p>, with no change to the existing run). A plain zero-length
# Range exactly at the paragraph-end offset does not anchor correctly, so we
# insert a one-character placeholder run, bookmark that, then delete the
# placeholder text -- leaving a correctly-anchored empty bookmark behind.
$target = $d.Content
[void]$target.Find.Execute("to which your program should make requests is:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$target.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $target)
$target.Text = ""
